# The author cleared out the stale "Turbosina " (AICM) figures that had
# been entered for the older, hidden history rows (years 2018-2022,
# sheet rows 33-85) - those cells go back to being blank while every
# other cell (formulas, formatting, the rest of the table) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G33:G85").ClearContents()
